$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.508.11"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "'2.644.00"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'604.80"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'154.68"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "'2.644.57"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +6.78%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'5.22"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "'28.05"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "'3.126.25"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'0.0000184"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "'67.465.36"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'2.643.99"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'11.33"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'365.01"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'7.65"
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'2.14"
$ws.Range("E23").Value = "  +6.49%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'10.07"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("D26").Value = "'66.12"
$ws.Range("E26").Value = "  -7.97%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("D28").Value = "'2.774.44"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "'583.53"
$ws.Range("E29").Value = "  -7.22%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'1.55"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Value = "'4.97"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'158.03"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").Value = "'19.51"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "'5.32"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("D44").Value = "'41.27"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'16.37"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").Value = "'156.77"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'21.06"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'0.630"
$ws.Range("E51").Value = "  +0.80%  "
